$d = $word.ActiveDocument

$pairs = @(
    @("59÷3=", "84÷4="),
    @("51÷4=", "22÷3="),
    @("47÷9=", "62÷4="),
    @("15÷2=", "12÷2="),
    @("79÷6=", "47÷7="),
    @("38÷7=", "97÷4="),
    @("74÷4=", "73÷6="),
    @("80÷4=", "83÷3="),
    @("33÷2=", "56÷6="),
    @("68÷5=", "69÷3="),
    @("40÷5=", "39÷2="),
    @("61÷7=", "52÷5="),
    @("47÷6=", "40÷4="),
    @("70÷5=", "43÷8="),
    @("18÷6=", "97÷9="),
    @("92÷8=", "79÷4="),
    @("80÷2=", "44÷7="),
    @("70÷9=", "12÷9="),
    @("16÷3=", "78÷7="),
    @("89÷9=", "47÷7="),
    @("91÷4=", "86÷8="),
    @("83÷9=", "91÷6="),
    @("81÷4=", "75÷3="),
    @("70÷7=", "46÷5="),
    @("17÷8=", "81÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
